$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow editing, re-protect afterward.
$ws.Unprotect()

# Update the confidential disclosure date (2021-04-06 -> 2021-04-08)
$ws.Range("A80").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) figures for each holding row
$ws.Range("D2").Value2 = 0.07632449312301398
$ws.Range("E2").Value2 = 0.01923377638780299
$ws.Range("D3").Value2 = 0.04696752703725051
$ws.Range("E3").Value2 = 0.006071251055836724
$ws.Range("D4").Value2 = 0.03698378206565611
$ws.Range("E4").Value2 = 0.01340536214485777
$ws.Range("D5").Value2 = 0.03349751314631211
$ws.Range("E5").Value2 = 0.006521639987230321
$ws.Range("D6").Value2 = 0.03106477286761779
$ws.Range("E6").Value2 = 0.00122636029174461
$ws.Range("D7").Value2 = 0.03206745829627309
$ws.Range("E7").Value2 = 0.005091490511516028
$ws.Range("D8").Value2 = 0.02968088269926756
$ws.Range("E8").Value2 = -0.00391174133610428
$ws.Range("D9").Value2 = 0.02894606700065963
$ws.Range("E9").Value2 = -0.0004172672786587617
$ws.Range("D10").Value2 = 0.02633745217455941
$ws.Range("E10").Value2 = 0.006655984222852185
$ws.Range("D11").Value2 = 0.02741752431367751
$ws.Range("E11").Value2 = 0.003802837501828416
$ws.Range("D12").Value2 = 0.02328074872143644
$ws.Range("E12").Value2 = -0.001279590531030217
$ws.Range("D13").Value2 = 0.02409108926640597
$ws.Range("E13").Value2 = -0.00873907615480618
$ws.Range("D14").Value2 = 0.02033098310386973
$ws.Range("E14").Value2 = 0.0185127786074557
$ws.Range("D15").Value2 = 0.01932901377679181
$ws.Range("E15").Value2 = 0.01434499110847653
$ws.Range("D16").Value2 = 0.02051349352589273
$ws.Range("E16").Value2 = 0.002704268881591698
$ws.Range("D17").Value2 = 0.01868738676345432
$ws.Range("E17").Value2 = -0.0006437768240344921
$ws.Range("D18").Value2 = 0.01780066205023033
$ws.Range("E18").Value2 = -0.006066522557701681
$ws.Range("D19").Value2 = 0.01485886903022517
$ws.Range("E19").Value2 = -0.01042587029510522
$ws.Range("D20").Value2 = 0.0135613884522128
$ws.Range("E20").Value2 = 0.008459271932325985
$ws.Range("D21").Value2 = 0.01644164447666795
$ws.Range("E21").Value2 = -0.0002235778849531966
$ws.Range("D22").Value2 = 0.01417231857793302
$ws.Range("E22").Value2 = 0.01794761237469022
$ws.Range("D23").Value2 = 0.01279100637526138
$ws.Range("E23").Value2 = 0.003627813234799993
$ws.Range("D24").Value2 = 0.01506133481620466
$ws.Range("E24").Value2 = -0.03006789524733278
$ws.Range("D25").Value2 = 0.01396455364028088
$ws.Range("E25").Value2 = -0.002051197899573309
$ws.Range("D26").Value2 = 0.01255984878608131
$ws.Range("E26").Value2 = 0.01415495955725876
$ws.Range("D27").Value2 = 0.01210880027253871
$ws.Range("E27").Value2 = 0.01048730484150773
$ws.Range("D28").Value2 = 0.01233484967046676
$ws.Range("E28").Value2 = 0.01207547169811307
$ws.Range("D29").Value2 = 0.01146908286340759
$ws.Range("E29").Value2 = -0.004945054945055039
$ws.Range("D30").Value2 = 0.01258009059066875
$ws.Range("E30").Value2 = 0.01399556756428555
$ws.Range("D31").Value2 = 0.01268645554496309
$ws.Range("E31").Value2 = 0.006321968841724868
$ws.Range("D32").Value2 = 0.0134256155931405
$ws.Range("E32").Value2 = 0.006272602169808073
$ws.Range("D33").Value2 = 0.01113251512203636
$ws.Range("E33").Value2 = 0.01099961404862992
$ws.Range("D34").Value2 = 0.01154833143797159
$ws.Range("E34").Value2 = -0.0237288135593221
$ws.Range("D35").Value2 = 0.009609653507568167
$ws.Range("E35").Value2 = 0.01912157026394601
$ws.Range("D36").Value2 = 0.01093744905235653
$ws.Range("E36").Value2 = -0.003003003003003046
$ws.Range("D37").Value2 = 0.01050043612973068
$ws.Range("E37").Value2 = -0.007501704932939335
$ws.Range("D38").Value2 = 0.01005263394333889
$ws.Range("E38").Value2 = -0.003419290497221783
$ws.Range("D39").Value2 = 0.009213076454012028
$ws.Range("E39").Value2 = 0.006259586286946117
$ws.Range("D40").Value2 = 0.009174502449043526
$ws.Range("E40").Value2 = -0.002747481475314228
$ws.Range("D41").Value2 = 0.009318964007254777
$ws.Range("E41").Value2 = 0.006526572473642744
$ws.Range("D42").Value2 = 0.009125807541781305
$ws.Range("E42").Value2 = -0.02094624285923541
$ws.Range("D43").Value2 = 0.009448339692235774
$ws.Range("E43").Value2 = -0.008276405675249787
$ws.Range("D44").Value2 = 0.009901106849564103
$ws.Range("E44").Value2 = 0.01633589847441619
$ws.Range("D45").Value2 = 0.008953274801736537
$ws.Range("E45").Value2 = -0.0112294845954507
$ws.Range("D46").Value2 = 0.009257474751809935
$ws.Range("E46").Value2 = 0.001392369813422611
$ws.Range("D47").Value2 = 0.008799456182914112
$ws.Range("E47").Value2 = -0.006770833333333171
$ws.Range("D48").Value2 = 0.006963180778076827
$ws.Range("E48").Value2 = 0.001590609916630292
$ws.Range("D49").Value2 = 0.008228580005422322
$ws.Range("E49").Value2 = 0
$ws.Range("D50").Value2 = 0.008105553754427477
$ws.Range("E50").Value2 = -0.001325205406837893
$ws.Range("D51").Value2 = 0.007699858340785961
$ws.Range("E51").Value2 = 0.0003038062584088674
$ws.Range("D52").Value2 = 0.007467936909923349
$ws.Range("E52").Value2 = -0.009397234528124465
$ws.Range("D53").Value2 = 0.00718965983695131
$ws.Range("E53").Value2 = -0.007768924302788749
$ws.Range("D54").Value2 = 0.007523840573064584
$ws.Range("E54").Value2 = 0.01598984771573608
$ws.Range("D55").Value2 = 0.006638261622364406
$ws.Range("E55").Value2 = 0.0006472491909386147
$ws.Range("D56").Value2 = 0.006567128865677437
$ws.Range("E56").Value2 = 0.006222739168362912
$ws.Range("D57").Value2 = 0.006645422638138263
$ws.Range("E57").Value2 = -0.009655172413793101
$ws.Range("D58").Value2 = 0.006246410839219013
$ws.Range("E58").Value2 = -0.004952538175815091
$ws.Range("D59").Value2 = 0.005440462383924111
$ws.Range("E59").Value2 = -0.005791505791505669
$ws.Range("D60").Value2 = 0.00660030823876297
$ws.Range("E60").Value2 = -0.007594662037539335
$ws.Range("D61").Value2 = 0.005326840933645599
$ws.Range("E61").Value2 = -0.01389137838322285
$ws.Range("D62").Value2 = 0.005752491711243593
$ws.Range("E62").Value2 = -0.002522905324658042
$ws.Range("D63").Value2 = 0.005230024000383067
$ws.Range("E63").Value2 = 0.004454505622900617
$ws.Range("D64").Value2 = 0.004740401481871961
$ws.Range("E64").Value2 = -0.00241701579116993
$ws.Range("D65").Value2 = 0.004659816184363503
$ws.Range("E65").Value2 = -0.004712728464879912
$ws.Range("D66").Value2 = 0.004398868769564194
$ws.Range("E66").Value2 = -0.006576805365631211
$ws.Range("D67").Value2 = 0.004337856915170942
$ws.Range("E67").Value2 = 0.008055995773903923
$ws.Range("D68").Value2 = 0.003660711263595129
$ws.Range("E68").Value2 = 0.03482003129890465
$ws.Range("D69").Value2 = 0.003997422225281833
$ws.Range("E69").Value2 = 0.0190964136003724
$ws.Range("D70").Value2 = 0.003680189226500017
$ws.Range("E70").Value2 = 0.001141552511415567
$ws.Range("D71").Value2 = 0.003162161345419282
$ws.Range("E71").Value2 = 0.01657683771909935
$ws.Range("D72").Value2 = 0.002700848709267486
$ws.Range("E72").Value2 = 0.01226711917135059
$ws.Range("D73").Value2 = 0.002611336012094288
$ws.Range("E73").Value2 = 0.01387593923106456
$ws.Range("D74").Value2 = 0.002355544528652155
$ws.Range("E74").Value2 = 0.01303175857805861
$ws.Range("D75").Value2 = 0.001820903090976075
$ws.Range("E75").Value2 = -0.01253211682659539
$ws.Range("D76").Value2 = 0.001938248269456993
$ws.Range("E76").Value2 = -0.01517241379310352
$ws.Range("E77").Value2 = 0.003487462421972953

# Restore sheet protection
$ws.Protect()

Write-Output "Update complete"
